# Add new player row (Carlton Carrington) at the top of the roster and
# refresh the rest of the table to match the latest roster snapshot,
# including team updates for De'Andre Hunter and Andrew Wiggins.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 'Carlton Carrington'
$ws.Range("B2").Value = 'PG,SG'
$ws.Range("C2").Value = 'Washington Wizards'

$ws.Range("A3").Value = 'Malik Beasley'
$ws.Range("B3").Value = 'SG,SF'
$ws.Range("C3").Value = 'Detroit Pistons'

$ws.Range("A4").Value = 'Bilal Coulibaly'
$ws.Range("B4").Value = 'SG,SF'
$ws.Range("C4").Value = 'Washington Wizards'

$ws.Range("A5").Value = 'Kelly Oubre Jr.'
$ws.Range("B5").Value = 'SG,SF'
$ws.Range("C5").Value = 'Philadelphia 76ers'

$ws.Range("A6").Value = 'De''Andre Hunter'
$ws.Range("B6").Value = 'SF,PF'
$ws.Range("C6").Value = 'Cleveland Cavaliers'

$ws.Range("A7").Value = 'Dyson Daniels'
$ws.Range("B7").Value = 'PG,SG,SF'
$ws.Range("C7").Value = 'Atlanta Hawks'

$ws.Range("A8").Value = 'Kel''el Ware'
$ws.Range("B8").Value = 'PF,C'
$ws.Range("C8").Value = 'Miami Heat'

$ws.Range("A9").Value = 'Victor Wembanyama'
$ws.Range("B9").Value = 'C'
$ws.Range("C9").Value = 'San Antonio Spurs'

$ws.Range("A10").Value = 'Donovan Mitchell'
$ws.Range("B10").Value = 'PG,SG'
$ws.Range("C10").Value = 'Cleveland Cavaliers'

$ws.Range("A11").Value = 'Jaden McDaniels'
$ws.Range("B11").Value = 'SF,PF'
$ws.Range("C11").Value = 'Minnesota Timberwolves'

$ws.Range("A12").Value = 'Alperen Sengün'
$ws.Range("B12").Value = 'C'
$ws.Range("C12").Value = 'Houston Rockets'

$ws.Range("A13").Value = 'Kristaps Porzingis'
$ws.Range("B13").Value = 'PF,C'
$ws.Range("C13").Value = 'Boston Celtics'

$ws.Range("A14").Value = 'Michael Porter Jr.'
$ws.Range("B14").Value = 'SF,PF'
$ws.Range("C14").Value = 'Denver Nuggets'

$ws.Range("A15").Value = 'Domantas Sabonis'
$ws.Range("B15").Value = 'C'
$ws.Range("C15").Value = 'Sacramento Kings'

$ws.Range("A16").Value = 'Josh Hart'
$ws.Range("B16").Value = 'SG,SF,PF'
$ws.Range("C16").Value = 'New York Knicks'

$ws.Range("A17").Value = 'Cam Thomas'
$ws.Range("B17").Value = 'SG,SF'
$ws.Range("C17").Value = 'Brooklyn Nets'

$ws.Range("A18").Value = 'Donte DiVincenzo'
$ws.Range("B18").Value = 'PG,SG,SF'
$ws.Range("C18").Value = 'Minnesota Timberwolves'

$ws.Range("A19").Value = 'Andrew Wiggins'
$ws.Range("B19").Value = 'SF,PF'
$ws.Range("C19").Value = 'Miami Heat'
